# Add new columns I (I0) and J (IF) to the active worksheet, matching the
# header style already used by the other header cells (e.g. H1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the existing header formatting (bold font, border, centered/top
# alignment) from H1 onto the two new header cells so the cell style matches
# the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for rows 2..29 (column I = I0, column J = IF)
$values = @{
    2  = @(9, 9)
    3  = @(7, 7)
    4  = @(7, 7)
    5  = @(7, 8)
    6  = @(9, 9)
    7  = @(11, 11)
    8  = @(11, 11)
    9  = @(7, 7)
    10 = @(7, 7)
    11 = @(5, 6)
    12 = @(4, 5)
    13 = @(8, 8)
    14 = @(5, 5)
    15 = @(6, 7)
    16 = @(6, 6)
    17 = @(7, 7)
    18 = @(6, 7)
    19 = @(9, 9)
    20 = @(6, 6)
    21 = @(7, 8)
    22 = @(9, 9)
    23 = @(6, 6)
    24 = @(9, 9)
    25 = @(7, 7)
    26 = @(4, 5)
    27 = @(9, 9)
    28 = @(8, 8)
    29 = @(8, 8)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
